# Update cryptos list cell values (prices + 1h volume %) and reorder a few rows
# per the authoritative diff (commit: "Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.226.42"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").Value = "3.135.56"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.27"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.83"
$ws.Range("E6").Value = "  -6.86%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.128.57"
$ws.Range("E8").Value = "  -3.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  -4.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -8.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.49"
$ws.Range("E11").Value = "  -5.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  -5.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  -8.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.60"
$ws.Range("E14").Value = "  -8.89%  "
$ws.Range("D15").Value = "3.665.32"
$ws.Range("E15").Value = "  -2.75%  "
$ws.Range("D16").Value = "64.202.06"
$ws.Range("E16").Value = "  -3.62%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "3.149.66"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("E19").Value = "  -5.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.33"
$ws.Range("E20").Value = "  -6.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("E21").Value = "  -5.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.705"
$ws.Range("E22").Value = "  -5.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  -5.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.52"
$ws.Range("E24").Value = "  -7.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.20"
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.86"
$ws.Range("E27").Value = "  -5.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.32"
$ws.Range("E28").Value = "  -8.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  -7.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.71"
$ws.Range("E30").Value = "  -2.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.113"
$ws.Range("E31").Value = "  -33.87%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.71"
$ws.Range("E33").Value = "  -6.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.92"
$ws.Range("E34").Value = "  -8.05%  "
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.10"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("E37").Value = "  -6.85%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "445.29"
$ws.Range("E38").Value = "  -10.14%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0703"
$ws.Range("E39").Value = "  -13.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.88"
$ws.Range("E40").Value = "  -11.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0392"
$ws.Range("E41").Value = "  -7.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  -8.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.35"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("D44").Value = "2.820.20"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.263"
$ws.Range("E45").Value = "  -9.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("E46").Value = "  -8.96%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.08"
$ws.Range("E48").Value = "  -7.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.28"
$ws.Range("E49").Value = "  -5.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.113"
$ws.Range("E50").Value = "  -5.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.64"
$ws.Range("E51").Value = "  -2.94%  "
